$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "create volunteering team for projects"
$ws.Range("B19").Value = "create volunteering team for projects"
$ws.Range("C19").Value = "CreateVolunteeringTeamForProjectsData.xlsx"
$ws.Range("D19").Value = "Yes"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "103"
$ws.Range("F19").Value = "TeamId"
